$wb = $excel.ActiveWorkbook

$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# --- Passenger (psgr) sheet: calibration tweak + header row height ---
$wsPsgr.Rows.Item(1).RowHeight = 45
$wsPsgr.Range("B2:H2").Value = 0.0755

# Select the whole header-to-data range on the psgr sheet (matches saved view state)
$wsPsgr.Activate()
$wsPsgr.Range("A1:H7").Select()

# --- Freight (frgt) sheet: calibration tweak (Reference scenario values) + header row height ---
$wsFrgt.Rows.Item(1).RowHeight = 45
$wsFrgt.Range("B3").Value = 0.0219
$wsFrgt.Range("C3").Value = 0.0219
$wsFrgt.Range("D3").Value = 0.081
$wsFrgt.Range("E3").Value = 0.0245
$wsFrgt.Range("F3").Value = 0.0219
$wsFrgt.Range("G3").Value = 0.0219
$wsFrgt.Range("H3").Value = 0.0219

# Freight sheet ends up the active/selected tab with D10 selected
$wsFrgt.Activate()
$wsFrgt.Range("D10").Select()
